$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store a value as literal text, matching the
# source workbook's inlineStr cells (prevents Excel from reinterpreting
# numeric-looking strings like '165.90' or '130.03' as numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

$ws.Range('D2').Value = '63.432.23'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.687.49'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '553.98'
$ws.Range('E5').Value = '  -3.83%  '
Set-TextValue $ws.Range('D6') '158.31'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -2.96%  '
$ws.Range('E9').Value = '  -3.99%  '
$ws.Range('E10').Value = '  -0.40%  '
Set-TextValue $ws.Range('D11') '0.368'
$ws.Range('E11').Value = '  -4.53%  '
Set-TextValue $ws.Range('D12') '5.38'
$ws.Range('E12').Value = '  -9.05%  '
$ws.Range('D13').Value = '3.164.12'
$ws.Range('E13').Value = '  -2.55%  '
Set-TextValue $ws.Range('D14') '26.34'
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('D15').Value = '63.282.80'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('E16').Value = '  -3.93%  '
$ws.Range('D17').Value = '2.690.32'
$ws.Range('E17').Value = '  -2.66%  '
Set-TextValue $ws.Range('D18') '12.03'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('E19').Value = '  -5.05%  '
Set-TextValue $ws.Range('D20') '341.94'
$ws.Range('E20').Value = '  -4.79%  '
$ws.Range('E21').Value = '  -4.84%  '
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('E23').Value = '  -3.95%  '
Set-TextValue $ws.Range('D24') '63.89'
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('E26').Value = '  +0.08%  '
Set-TextValue $ws.Range('D27') '8.18'
$ws.Range('E27').Value = '  -4.24%  '
$ws.Range('D28').Value = '0.0₃0854'
$ws.Range('E28').Value = '  -5.37%  '
$ws.Range('E29').Value = '  -0.77%  '
Set-TextValue $ws.Range('D30') '1.33'
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('E31').Value = '  -4.63%  '
Set-TextValue $ws.Range('D32') '165.90'
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -3.08%  '
Set-TextValue $ws.Range('D35') '19.54'
$ws.Range('E35').Value = '  -3.31%  '
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('E37').Value = '  -2.14%  '
Set-TextValue $ws.Range('D38') '339.50'
$ws.Range('E38').Value = '  -2.20%  '
Set-TextValue $ws.Range('D39') '0.946'
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('E40').Value = '  -4.14%  '
Set-TextValue $ws.Range('D41') '38.14'
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  -6.06%  '
Set-TextValue $ws.Range('D43') '20.78'
$ws.Range('E43').Value = '  -5.21%  '
Set-TextValue $ws.Range('D44') '20.27'
Set-TextValue $ws.Range('D45') '0.619'
$ws.Range('E45').Value = '  -1.38%  '
Set-TextValue $ws.Range('D46') '0.0562'
$ws.Range('E46').Value = '  -4.35%  '
$ws.Range('E47').Value = '  +0.02%  '
Set-TextValue $ws.Range('D48') '11.08'
$ws.Range('E48').Value = '  +0.30%  '
Set-TextValue $ws.Range('D49') '130.03'
$ws.Range('E49').Value = '  -5.31%  '
Set-TextValue $ws.Range('D50') '0.0971'
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.094.63'
$ws.Range('E51').Value = '  -1.59%  '
